# Update Leve profit-calc columns (H-N) across all class sheets
# per the scheduled-runner recalculation pass.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H92").Value = 20585.312
$ws.Range("I92").Value = 1545.1515
$ws.Range("K92").Value = 1545.1515
$ws.Range("M92").Value = -297.1514999999999

$ws.Range("H96").Value = 1199.2858
$ws.Range("I96").Value = 1239.25
$ws.Range("J96").Value = 1146
$ws.Range("K96").Value = 3717.75
$ws.Range("L96").Value = 3438
$ws.Range("M96").Value = -2344.75
$ws.Range("N96").Value = -6184

$ws.Range("H132").Value = 3075.6924
$ws.Range("I132").Value = 3187.2
$ws.Range("K132").Value = 9561.599999999999
$ws.Range("M132").Value = -7031.599999999999

$ws.Range("H137").Value = 1106.8572
$ws.Range("J137").Value = 1207.6364
$ws.Range("L137").Value = 3622.9092
$ws.Range("N137").Value = -8722.9092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 909.625
$ws.Range("I88").Value = 857.6667
$ws.Range("J88").Value = 940.8
$ws.Range("K88").Value = 857.6667
$ws.Range("L88").Value = 940.8
$ws.Range("M88").Value = -451.6667
$ws.Range("N88").Value = -1752.8

$ws.Range("H91").Value = 909.625
$ws.Range("I91").Value = 857.6667
$ws.Range("J91").Value = 940.8
$ws.Range("K91").Value = 857.6667
$ws.Range("L91").Value = 940.8
$ws.Range("M91").Value = 546.3333
$ws.Range("N91").Value = -3748.8

$ws.Range("H108").Value = 64346.668
$ws.Range("I108").Value = 37660
$ws.Range("J108").Value = 69684
$ws.Range("K108").Value = 37660
$ws.Range("L108").Value = 69684
$ws.Range("M108").Value = -33820
$ws.Range("N108").Value = -77364

$ws.Range("H110").Value = 1316.5555
$ws.Range("I110").Value = 1104.5
$ws.Range("K110").Value = 1104.5
$ws.Range("M110").Value = 940.5

$ws.Range("H122").Value = 61106.465
$ws.Range("I122").Value = 924.75
$ws.Range("K122").Value = 2774.25
$ws.Range("M122").Value = -324.25

$ws.Range("H133").Value = 88833.336
$ws.Range("J133").Value = 88833.336
$ws.Range("L133").Value = 88833.336
$ws.Range("N133").Value = -93893.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2887.15
$ws.Range("I134").Value = 2887.15
$ws.Range("K134").Value = 8661.450000000001
$ws.Range("M134").Value = -6126.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1587.625
$ws.Range("I22").Value = 1450.25
$ws.Range("J22").Value = 1725
$ws.Range("K22").Value = 1450.25
$ws.Range("L22").Value = 1725
$ws.Range("M22").Value = -1100.25
$ws.Range("N22").Value = -2425

$ws.Range("H43").Value = 4183.1665
$ws.Range("J43").Value = 4183.1665
$ws.Range("L43").Value = 4183.1665
$ws.Range("N43").Value = -4551.1665

$ws.Range("H58").Value = 75761380
$ws.Range("I58").Value = 111116056
$ws.Range("J58").Value = 62503376
$ws.Range("K58").Value = 111116056
$ws.Range("L58").Value = 62503376
$ws.Range("M58").Value = -111115853
$ws.Range("N58").Value = -62503782

$ws.Range("H62").Value = 63134.723
$ws.Range("I62").Value = 91479.164
$ws.Range("J62").Value = 6445.8335
$ws.Range("K62").Value = 91479.164
$ws.Range("L62").Value = 6445.8335
$ws.Range("M62").Value = -90855.164
$ws.Range("N62").Value = -7693.8335

$ws.Range("H65").Value = 63134.723
$ws.Range("I65").Value = 91479.164
$ws.Range("J65").Value = 6445.8335
$ws.Range("K65").Value = 457395.82
$ws.Range("L65").Value = 32229.1675
$ws.Range("M65").Value = -454275.82
$ws.Range("N65").Value = -38469.1675

$ws.Range("H68").Value = 55627
$ws.Range("I68").Value = 40475
$ws.Range("J68").Value = 64285.285
$ws.Range("K68").Value = 40475
$ws.Range("L68").Value = 64285.285
$ws.Range("M68").Value = -39726
$ws.Range("N68").Value = -65783.285

$ws.Range("H71").Value = 55627
$ws.Range("I71").Value = 40475
$ws.Range("J71").Value = 64285.285
$ws.Range("K71").Value = 121425
$ws.Range("L71").Value = 192855.855
$ws.Range("M71").Value = -117681
$ws.Range("N71").Value = -200343.855

$ws.Range("H101").Value = 4183.1665
$ws.Range("J101").Value = 4183.1665
$ws.Range("L101").Value = 4183.1665
$ws.Range("N101").Value = -10673.1665

$ws.Range("H132").Value = 3876.318
$ws.Range("I132").Value = 2646.3684
$ws.Range("J132").Value = 11666
$ws.Range("K132").Value = 7939.1052
$ws.Range("L132").Value = 34998
$ws.Range("M132").Value = -5409.1052
$ws.Range("N132").Value = -40058

$ws.Range("H134").Value = 3055.5
$ws.Range("I134").Value = 3141.4546
$ws.Range("K134").Value = 9424.363799999999
$ws.Range("M134").Value = -6889.363799999999

$ws.Range("H136").Value = 75761380
$ws.Range("I136").Value = 111116056
$ws.Range("J136").Value = 62503376
$ws.Range("K136").Value = 333348168
$ws.Range("L136").Value = 187510128
$ws.Range("M136").Value = -333345618
$ws.Range("N136").Value = -187515228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1915.4286
$ws.Range("I59").Value = 1496.4
$ws.Range("J59").Value = 2963
$ws.Range("K59").Value = 4489.200000000001
$ws.Range("L59").Value = 8889
$ws.Range("M59").Value = -3949.200000000001
$ws.Range("N59").Value = -9969

$ws.Range("H61").Value = 2022.375
$ws.Range("I61").Value = 1032.7142
$ws.Range("J61").Value = 8950
$ws.Range("K61").Value = 3098.1426
$ws.Range("L61").Value = 26850
$ws.Range("M61").Value = -2883.1426
$ws.Range("N61").Value = -27280

$ws.Range("H69").Value = 5626.4443
$ws.Range("I69").Value = 935
$ws.Range("J69").Value = 15009.333
$ws.Range("K69").Value = 2805
$ws.Range("L69").Value = 45027.999
$ws.Range("M69").Value = -1994
$ws.Range("N69").Value = -46649.999

$ws.Range("H72").Value = 5626.4443
$ws.Range("I72").Value = 935
$ws.Range("J72").Value = 15009.333
$ws.Range("K72").Value = 8415
$ws.Range("L72").Value = 135083.997
$ws.Range("M72").Value = -4359
$ws.Range("N72").Value = -143195.997

$ws.Range("H113").Value = 52633948
$ws.Range("J113").Value = 76924890
$ws.Range("L113").Value = 230774670
$ws.Range("N113").Value = -230779010

$ws.Range("H139").Value = 5146.8823
$ws.Range("I139").Value = 3167.889
$ws.Range("K139").Value = 9503.667000000001
$ws.Range("M139").Value = -4363.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 53947.184
$ws.Range("I102").Value = 114294.89
$ws.Range("J102").Value = 12168
$ws.Range("K102").Value = 114294.89
$ws.Range("L102").Value = 12168
$ws.Range("M102").Value = -112672.89
$ws.Range("N102").Value = -15412

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 8601.200000000001
$ws.Range("J32").Value = 9998.25
$ws.Range("L32").Value = 9998.25
$ws.Range("N32").Value = -10632.25

$ws.Range("H40").Value = 5061.6
$ws.Range("I40").Value = 4920.6924
$ws.Range("K40").Value = 4920.6924
$ws.Range("M40").Value = -4784.6924

$ws.Range("H122").Value = 3320.35
$ws.Range("I122").Value = 3320.35
$ws.Range("K122").Value = 9961.049999999999
$ws.Range("M122").Value = -7511.049999999999

$ws.Range("H132").Value = 2747.7646
$ws.Range("I132").Value = 2428.9062
$ws.Range("K132").Value = 7286.7186
$ws.Range("M132").Value = -4756.7186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7000
$ws.Range("I54").Value = 7000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 7000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -6480
$ws.Range("N54").ClearContents()

$ws.Range("H62").Value = 9934.571
$ws.Range("J62").Value = 10509.223
$ws.Range("L62").Value = 10509.223
$ws.Range("N62").Value = -11757.223

$ws.Range("H65").Value = 9934.571
$ws.Range("J65").Value = 10509.223
$ws.Range("L65").Value = 52546.115
$ws.Range("N65").Value = -58786.115

$ws.Range("H74").Value = 11595
$ws.Range("I74").Value = 14000
$ws.Range("J74").Value = 10993.75
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 10993.75
$ws.Range("M74").Value = -13064
$ws.Range("N74").Value = -12865.75

$ws.Range("H77").Value = 11595
$ws.Range("I77").Value = 14000
$ws.Range("J77").Value = 10993.75
$ws.Range("K77").Value = 42000
$ws.Range("L77").Value = 32981.25
$ws.Range("M77").Value = -37320
$ws.Range("N77").Value = -42341.25

$ws.Range("H109").Value = 53227.816
$ws.Range("J109").Value = 53227.816
$ws.Range("L109").Value = 53227.816
$ws.Range("N109").Value = -56001.816

$ws.Range("H122").Value = 2386.9412
$ws.Range("I122").Value = 1827.4286
$ws.Range("K122").Value = 5482.2858
$ws.Range("M122").Value = -3032.2858

$ws.Range("H132").Value = 3860.8
$ws.Range("I132").Value = 3849.75
$ws.Range("J132").Value = 3905
$ws.Range("K132").Value = 11549.25
$ws.Range("L132").Value = 11715
$ws.Range("M132").Value = -9019.25
$ws.Range("N132").Value = -16775
